# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header style used by the other header cells (B1:G1), and populate the
# value for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it picks up the same bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# New data cell for row 2 - plain numeric value, same as the diff (no
# special style applied, just like B2:G2 data cells other than A2).
$ws.Range("H2").Value = 1
